# Natmi following Dr Hou advice
# Update NATMI LR-pair (Fn1-Itgb3) statistics for rows 2-17 (columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T)
# to reflect the recomputed expression/specificity values from the revised analysis
# (ligand/receptor-expressing cell counts now 3 instead of 1, plus corresponding
# average/total expression and derived-specificity recalculations).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 19.95578266666667
$ws.Cells.Item(2, 8).Value = 59.867348
$ws.Cells.Item(2, 9).Value = 0.0117373419656925
$ws.Cells.Item(2, 10).Value = 0.0117373419656925
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 3.778439
$ws.Cells.Item(2, 14).Value = 11.335317
$ws.Cells.Item(2, 15).Value = 0.4252971528324392
$ws.Cells.Item(2, 16).Value = 0.4252971528324392
$ws.Cells.Item(2, 17).Value = 75.40170750325734
$ws.Cells.Item(2, 18).Value = 678.615367529316
$ws.Cells.Item(2, 19).Value = 0.004991858119829724
$ws.Cells.Item(2, 20).Value = 0.004991858119829724

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 19.95578266666667
$ws.Cells.Item(3, 8).Value = 59.867348
$ws.Cells.Item(3, 9).Value = 0.0117373419656925
$ws.Cells.Item(3, 10).Value = 0.0117373419656925
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 4.333403333333333
$ws.Cells.Item(3, 14).Value = 13.00021
$ws.Cells.Item(3, 15).Value = 0.4877633593505858
$ws.Cells.Item(3, 16).Value = 0.4877633593505858
$ws.Cells.Item(3, 17).Value = 86.47645512700889
$ws.Cells.Item(3, 18).Value = 778.28809614308
$ws.Cells.Item(3, 19).Value = 0.00572504534703278
$ws.Cells.Item(3, 20).Value = 0.005725045347032781

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 19.95578266666667
$ws.Cells.Item(4, 8).Value = 59.867348
$ws.Cells.Item(4, 9).Value = 0.0117373419656925
$ws.Cells.Item(4, 10).Value = 0.0117373419656925
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 0.2909853333333334
$ws.Cells.Item(4, 14).Value = 0.8729560000000001
$ws.Cells.Item(4, 15).Value = 0.03275300561492853
$ws.Cells.Item(4, 16).Value = 0.03275300561492853
$ws.Cells.Item(4, 17).Value = 5.806840071187557
$ws.Cells.Item(4, 18).Value = 52.261560640688
$ws.Cells.Item(4, 19).Value = 0.0003844332273066626
$ws.Cells.Item(4, 20).Value = 0.0003844332273066626

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 19.95578266666667
$ws.Cells.Item(5, 8).Value = 59.867348
$ws.Cells.Item(5, 9).Value = 0.0117373419656925
$ws.Cells.Item(5, 10).Value = 0.0117373419656925
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 0.4814053333333333
$ws.Cells.Item(5, 14).Value = 1.444216
$ws.Cells.Item(5, 15).Value = 0.0541864822020464
$ws.Cells.Item(5, 16).Value = 0.05418648220204641
$ws.Cells.Item(5, 17).Value = 9.606820206574222
$ws.Cells.Item(5, 18).Value = 86.461381859168
$ws.Cells.Item(5, 19).Value = 0.0006360052715233288
$ws.Cells.Item(5, 20).Value = 0.0006360052715233289

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 1637.343343333333
$ws.Cells.Item(6, 8).Value = 4912.03003
$ws.Cells.Item(6, 9).Value = 0.9630320723052701
$ws.Cells.Item(6, 10).Value = 0.9630320723052702
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 3.778439
$ws.Cells.Item(6, 14).Value = 11.335317
$ws.Cells.Item(6, 15).Value = 0.4252971528324392
$ws.Cells.Item(6, 16).Value = 0.4252971528324392
$ws.Cells.Item(6, 17).Value = 6186.601944841056
$ws.Cells.Item(6, 18).Value = 55679.41750356951
$ws.Cells.Item(6, 19).Value = 0.4095747984377551
$ws.Cells.Item(6, 20).Value = 0.4095747984377552

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 1637.343343333333
$ws.Cells.Item(7, 8).Value = 4912.03003
$ws.Cells.Item(7, 9).Value = 0.9630320723052701
$ws.Cells.Item(7, 10).Value = 0.9630320723052702
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 4.333403333333333
$ws.Cells.Item(7, 14).Value = 13.00021
$ws.Cells.Item(7, 15).Value = 0.4877633593505858
$ws.Cells.Item(7, 16).Value = 0.4877633593505858
$ws.Cells.Item(7, 17).Value = 7095.26910181181
$ws.Cells.Item(7, 18).Value = 63857.4219163063
$ws.Cells.Item(7, 19).Value = 0.4697317587499748
$ws.Cells.Item(7, 20).Value = 0.4697317587499749

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 1637.343343333333
$ws.Cells.Item(8, 8).Value = 4912.03003
$ws.Cells.Item(8, 9).Value = 0.9630320723052701
$ws.Cells.Item(8, 10).Value = 0.9630320723052702
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 0.2909853333333334
$ws.Cells.Item(8, 14).Value = 0.8729560000000001
$ws.Cells.Item(8, 15).Value = 0.03275300561492853
$ws.Cells.Item(8, 16).Value = 0.03275300561492853
$ws.Cells.Item(8, 17).Value = 476.4428985409645
$ws.Cells.Item(8, 18).Value = 4287.98608686868
$ws.Cells.Item(8, 19).Value = 0.03154219487157078
$ws.Cells.Item(8, 20).Value = 0.03154219487157078

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 1637.343343333333
$ws.Cells.Item(9, 8).Value = 4912.03003
$ws.Cells.Item(9, 9).Value = 0.9630320723052701
$ws.Cells.Item(9, 10).Value = 0.9630320723052702
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 0.4814053333333333
$ws.Cells.Item(9, 14).Value = 1.444216
$ws.Cells.Item(9, 15).Value = 0.0541864822020464
$ws.Cells.Item(9, 16).Value = 0.05418648220204641
$ws.Cells.Item(9, 17).Value = 788.2258179784976
$ws.Cells.Item(9, 18).Value = 7094.032361806479
$ws.Cells.Item(9, 19).Value = 0.05218332024596938
$ws.Cells.Item(9, 20).Value = 0.0521833202459694

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 17.50081933333334
$ws.Cells.Item(10, 8).Value = 52.502458
$ws.Cells.Item(10, 9).Value = 0.01029341242216722
$ws.Cells.Item(10, 10).Value = 0.01029341242216722
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 3.778439
$ws.Cells.Item(10, 14).Value = 11.335317
$ws.Cells.Item(10, 15).Value = 0.4252971528324392
$ws.Cells.Item(10, 16).Value = 0.4252971528324392
$ws.Cells.Item(10, 17).Value = 66.12577830102067
$ws.Cells.Item(10, 18).Value = 595.132004709186
$ws.Cells.Item(10, 19).Value = 0.00437775899607778
$ws.Cells.Item(10, 20).Value = 0.004377758996077781

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 17.50081933333334
$ws.Cells.Item(11, 8).Value = 52.502458
$ws.Cells.Item(11, 9).Value = 0.01029341242216722
$ws.Cells.Item(11, 10).Value = 0.01029341242216722
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 4.333403333333333
$ws.Cells.Item(11, 14).Value = 13.00021
$ws.Cells.Item(11, 15).Value = 0.4877633593505858
$ws.Cells.Item(11, 16).Value = 0.4877633593505858
$ws.Cells.Item(11, 17).Value = 75.83810883513112
$ws.Cells.Item(11, 18).Value = 682.54297951618
$ws.Cells.Item(11, 19).Value = 0.005020749422217332
$ws.Cells.Item(11, 20).Value = 0.005020749422217334

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 17.50081933333334
$ws.Cells.Item(12, 8).Value = 52.502458
$ws.Cells.Item(12, 9).Value = 0.01029341242216722
$ws.Cells.Item(12, 10).Value = 0.01029341242216722
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 0.2909853333333334
$ws.Cells.Item(12, 14).Value = 0.8729560000000001
$ws.Cells.Item(12, 15).Value = 0.03275300561492853
$ws.Cells.Item(12, 16).Value = 0.03275300561492853
$ws.Cells.Item(12, 17).Value = 5.092481747316446
$ws.Cells.Item(12, 18).Value = 45.832335725848
$ws.Cells.Item(12, 19).Value = 0.000337140194860018
$ws.Cells.Item(12, 20).Value = 0.0003371401948600181

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 17.50081933333334
$ws.Cells.Item(13, 8).Value = 52.502458
$ws.Cells.Item(13, 9).Value = 0.01029341242216722
$ws.Cells.Item(13, 10).Value = 0.01029341242216722
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 0.4814053333333333
$ws.Cells.Item(13, 14).Value = 1.444216
$ws.Cells.Item(13, 15).Value = 0.0541864822020464
$ws.Cells.Item(13, 16).Value = 0.05418648220204641
$ws.Cells.Item(13, 17).Value = 8.424987764769778
$ws.Cells.Item(13, 18).Value = 75.82488988292801
$ws.Cells.Item(13, 19).Value = 0.0005577638090120873
$ws.Cells.Item(13, 20).Value = 0.0005577638090120874

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 25.39612333333333
$ws.Cells.Item(14, 8).Value = 76.18836999999999
$ws.Cells.Item(14, 9).Value = 0.01493717330687017
$ws.Cells.Item(14, 10).Value = 0.01493717330687017
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 3.778439
$ws.Cells.Item(14, 14).Value = 11.335317
$ws.Cells.Item(14, 15).Value = 0.4252971528324392
$ws.Cells.Item(14, 16).Value = 0.4252971528324392
$ws.Cells.Item(14, 17).Value = 95.95770285147667
$ws.Cells.Item(14, 18).Value = 863.6193256632899
$ws.Cells.Item(14, 19).Value = 0.006352737278776594
$ws.Cells.Item(14, 20).Value = 0.006352737278776594

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 25.39612333333333
$ws.Cells.Item(15, 8).Value = 76.18836999999999
$ws.Cells.Item(15, 9).Value = 0.01493717330687017
$ws.Cells.Item(15, 10).Value = 0.01493717330687017
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 4.333403333333333
$ws.Cells.Item(15, 14).Value = 13.00021
$ws.Cells.Item(15, 15).Value = 0.4877633593505858
$ws.Cells.Item(15, 16).Value = 0.4877633593505858
$ws.Cells.Item(15, 17).Value = 110.0516455064111
$ws.Cells.Item(15, 18).Value = 990.4648095576998
$ws.Cells.Item(15, 19).Value = 0.007285805831360892
$ws.Cells.Item(15, 20).Value = 0.007285805831360893

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 25.39612333333333
$ws.Cells.Item(16, 8).Value = 76.18836999999999
$ws.Cells.Item(16, 9).Value = 0.01493717330687017
$ws.Cells.Item(16, 10).Value = 0.01493717330687017
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 0.2909853333333334
$ws.Cells.Item(16, 14).Value = 0.8729560000000001
$ws.Cells.Item(16, 15).Value = 0.03275300561492853
$ws.Cells.Item(16, 16).Value = 0.03275300561492853
$ws.Cells.Item(16, 17).Value = 7.389899413524445
$ws.Cells.Item(16, 18).Value = 66.50909472172
$ws.Cells.Item(16, 19).Value = 0.0004892373211910793
$ws.Cells.Item(16, 20).Value = 0.0004892373211910793

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 25.39612333333333
$ws.Cells.Item(17, 8).Value = 76.18836999999999
$ws.Cells.Item(17, 9).Value = 0.01493717330687017
$ws.Cells.Item(17, 10).Value = 0.01493717330687017
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 0.4814053333333333
$ws.Cells.Item(17, 14).Value = 1.444216
$ws.Cells.Item(17, 15).Value = 0.0541864822020464
$ws.Cells.Item(17, 16).Value = 0.0541864822020464
$ws.Cells.Item(17, 17).Value = 12.22582921865778
$ws.Cells.Item(17, 18).Value = 110.03246296792
$ws.Cells.Item(17, 19).Value = 0.0008093928755416029
$ws.Cells.Item(17, 20).Value = 0.000809392875541603
